# Append one new data row (row 87) to Sheet1, mirroring the existing
# daily log rows: date text, weekday-kanji text, an hour number, and the
# ranking number. Also bumps the sheet's used-range dimension (handled
# automatically by Excel once the new row's cells are populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# Column A / B hold human-typed-looking strings ("2025/10/10", "金").
# A plain .Value assignment would let Excel's smart-type detection turn
# the date-shaped string into a real date serial (like pasting into a
# General-formatted cell). Temporarily forcing Text number format keeps
# the literal string, then reverting the format to General (and the
# style back to Normal) afterwards leaves the cell indistinguishable
# from the sheet's other untouched text cells.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/10"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "金"

# Column C / D are plain numbers, same as every other row.
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 201
